# Update the bird records table (rows 2-7 -> rows 2-11), per the
# "checking name or password is valid, adding some message boxs" commit:
# several subspecies/date/cage values were corrected and 4 new bird
# records were appended, then the table was re-sorted by SerialNumber.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure rows 8-11 inherit the same formatting (incl. the D-column
# date number format) as the existing data rows before we write values
# into them.
$ws.Range("A2:H2").Copy()
$ws.Range("A8:H11").PasteSpecial(-4122)

# SerialNumber, Strain, SubSpecies, DateOfBird, Gender, CageNumber,
# FatherSerialNumber, MotherSerialNumber
$data = @(
    @(2,  "Golden Australian ", "Coastal cities",  45077, "Male",   2, 3,  3),
    @(3,  "Golden European",    "East Europe",     45077, "Male",   1, 3,  1),
    @(4,  "Golden European",    "West Europe",     45077, "Male",   4, 1,  1),
    @(5,  "Golden Australian ", "Coastal cities",  45047, "Male",   3, 10, 10),
    @(6,  "Golden European",    "East Europe",     45077, "Female", 4, 10, 10),
    @(7,  "Golden European",    "West Europe",     45077, "Male",   2, 3,  3),
    @(8,  "Golden European",    "West Europe",     45054, "Female", 1, 2,  2),
    @(9,  "Golden Amrican",     "Central America", 45047, "Male",   2, 4,  4),
    @(10, "Golden European",    "East Europe",     45077, "Male",   1, 2,  2),
    @(11, "Golden Australian ", "Coastal cities",  45047, "Male",   2, 3,  3)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Refresh the sorted-range marker to cover the now-larger table.
$so = $ws.Sort
$so.SetRange($ws.Range("A2:H11"))
$so.Header = -4142
$so.Apply()
